$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4569.6
$ws.Range("I86").Value = 3235.2942
$ws.Range("J86").Value = 7405
$ws.Range("K86").Value = 3235.2942
$ws.Range("L86").Value = 7405
$ws.Range("M86").Value = -2112.2942
$ws.Range("N86").Value = -9651

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4569.6
$ws.Range("I89").Value = 3235.2942
$ws.Range("J89").Value = 7405
$ws.Range("K89").Value = 16176.471
$ws.Range("L89").Value = 37025
$ws.Range("M89").Value = -10560.471
$ws.Range("N89").Value = -48257

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9371.083000000001
$ws.Range("J116").Value = 6496.5
$ws.Range("L116").Value = 6496.5
$ws.Range("N116").Value = -13380.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 79600
$ws.Range("J136").Value = 79600
$ws.Range("L136").Value = 79600
$ws.Range("N136").Value = -89800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 185000
$ws.Range("J139").Value = 250000
$ws.Range("L139").Value = 250000
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 52999.316
$ws.Range("I5").Value = 91180.82000000001
$ws.Range("K5").Value = 91180.82000000001
$ws.Range("M5").Value = -91068.82000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8681.666999999999
$ws.Range("I45").Value = 13419.777
$ws.Range("K45").Value = 13419.777
$ws.Range("M45").Value = -13042.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 27999.5
$ws.Range("J62").Value = 27999.5
$ws.Range("L62").Value = 27999.5
$ws.Range("N62").Value = -29247.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 27999.5
$ws.Range("J65").Value = 27999.5
$ws.Range("L65").Value = 83998.5
$ws.Range("N65").Value = -90238.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2084.9092
$ws.Range("I122").Value = 1565.375
$ws.Range("K122").Value = 4696.125
$ws.Range("M122").Value = -2246.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 52999.316
$ws.Range("I4").Value = 91180.82000000001
$ws.Range("K4").Value = 91180.82000000001
$ws.Range("M4").Value = -91065.82000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 203571.19
$ws.Range("I22").Value = 314.76923
$ws.Range("J22").Value = 342641.38
$ws.Range("K22").Value = 314.76923
$ws.Range("L22").Value = 342641.38
$ws.Range("M22").Value = -141.76923
$ws.Range("N22").Value = -342987.38

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 19999
$ws.Range("I138").Value = 19999
$ws.Range("K138").Value = 19999
$ws.Range("M138").Value = -14859

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 645.8182
$ws.Range("I7").Value = 807.4286
$ws.Range("J7").Value = 363
$ws.Range("K7").Value = 807.4286
$ws.Range("L7").Value = 363
$ws.Range("M7").Value = -694.4286
$ws.Range("N7").Value = -589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 530.7143
$ws.Range("J22").Value = 623
$ws.Range("L22").Value = 623
$ws.Range("N22").Value = -1323

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32827.47
$ws.Range("I31").Value = 45048.22
$ws.Range("J31").Value = 7275
$ws.Range("K31").Value = 45048.22
$ws.Range("L31").Value = 7275
$ws.Range("M31").Value = -44753.22
$ws.Range("N31").Value = -7865

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 32827.47
$ws.Range("I34").Value = 45048.22
$ws.Range("J34").Value = 7275
$ws.Range("K34").Value = 45048.22
$ws.Range("L34").Value = 7275
$ws.Range("M34").Value = -44846.22
$ws.Range("N34").Value = -7679

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3722.2144
$ws.Range("I99").Value = 3509.25
$ws.Range("K99").Value = 3509.25
$ws.Range("M99").Value = -2011.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1425.2727
$ws.Range("I122").Value = 1328
$ws.Range("K122").Value = 3984
$ws.Range("M122").Value = -1534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3722.2144
$ws.Range("I126").Value = 3509.25
$ws.Range("K126").Value = 10527.75
$ws.Range("M126").Value = -8057.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 74316.664
$ws.Range("J135").Value = 74316.664
$ws.Range("L135").Value = 74316.664
$ws.Range("N135").Value = -84456.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 683.3333
$ws.Range("J12").Value = 1362.3334
$ws.Range("L12").Value = 4087.0002
$ws.Range("N12").Value = -4433.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5788
$ws.Range("J75").Value = 5261.5
$ws.Range("L75").Value = 15784.5
$ws.Range("N75").Value = -17780.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 5788
$ws.Range("J78").Value = 5261.5
$ws.Range("L78").Value = 47353.5
$ws.Range("N78").Value = -57337.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8225
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 9633.333000000001
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 9633.333000000001
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -11629.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8225
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 9633.333000000001
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 48166.665
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -58150.665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2233.1333
$ws.Range("I102").Value = 2290.3333
$ws.Range("K102").Value = 2290.3333
$ws.Range("M102").Value = -668.3332999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2229.2144
$ws.Range("I122").Value = 2158.4546
$ws.Range("J122").Value = 2488.6667
$ws.Range("K122").Value = 6475.3638
$ws.Range("L122").Value = 7466.000100000001
$ws.Range("M122").Value = -4025.3638
$ws.Range("N122").Value = -12366.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 12753998
$ws.Range("J23").Value = 25004996
$ws.Range("L23").Value = 25004996
$ws.Range("N23").Value = -25005456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1560.4615
$ws.Range("J46").Value = 1499.6
$ws.Range("L46").Value = 1499.6
$ws.Range("N46").Value = -1875.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 38212
$ws.Range("J76").Value = 38212
$ws.Range("L76").Value = 38212
$ws.Range("N76").Value = -38888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 38212
$ws.Range("J79").Value = 38212
$ws.Range("L79").Value = 38212
$ws.Range("N79").Value = -40552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 48000
$ws.Range("J88").Value = 48000
$ws.Range("L88").Value = 48000
$ws.Range("N88").Value = -48856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 48000
$ws.Range("J91").Value = 48000
$ws.Range("L91").Value = 48000
$ws.Range("N91").Value = -50964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4348.9546
$ws.Range("I122").Value = 3908.0715
$ws.Range("J122").Value = 5120.5
$ws.Range("K122").Value = 11724.2145
$ws.Range("L122").Value = 15361.5
$ws.Range("M122").Value = -9274.2145
$ws.Range("N122").Value = -20261.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4544.5
$ws.Range("I4").Value = 4816.3335
$ws.Range("J4").Value = 4428
$ws.Range("K4").Value = 4816.3335
$ws.Range("L4").Value = 4428
$ws.Range("M4").Value = -4703.3335
$ws.Range("N4").Value = -4654

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 42246.125
$ws.Range("J69").Value = 47487.715
$ws.Range("L69").Value = 47487.715
$ws.Range("N69").Value = -48985.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 42246.125
$ws.Range("J72").Value = 47487.715
$ws.Range("L72").Value = 142463.145
$ws.Range("N72").Value = -149951.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 83293.336
$ws.Range("J95").Value = 83293.336
$ws.Range("L95").Value = 83293.336
$ws.Range("N95").Value = -88785.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2824.2727
$ws.Range("I122").Value = 2771.3447
$ws.Range("J122").Value = 2926.6
$ws.Range("K122").Value = 8314.034100000001
$ws.Range("L122").Value = 8779.799999999999
$ws.Range("M122").Value = -5864.034100000001
